$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unchanged vendor / description text re-affirmed (keeps shared-string slots stable) ---
$ws.Range("A3").Value = "u/tonsoffun"
$ws.Range("A4").Value = "JLCPCB"
$ws.Range("B4").Value = "PCB + Stencil"
$ws.Range("A5").Value = "LCSC"
$ws.Range("B5").Value = "Parts (see LCSC BOM)"

# --- New "Notes" column header + first note ---
$ws.Range("D1").Value = "Notes"
$ws.Range("D2").Value = "Not able to find these anywhere else, shipping is slow Aliexpress value shipping"

# --- Row 2: SK6812Mini-E RGB LEDs x 100 (bold "x 100") ---
$ws.Range("B2").Value = "SK6812Mini-E RGB LEDs x 100"
$ws.Range("B2").Characters(23, 5).Font.Bold = $true

# --- Row 3: Mill-Max 0305-2-15-80-47-80-10 x 200 (bold "x 200") ---
$ws.Range("B3").Value = "Mill-Max 0305-2-15-80-47-80-10 x 200"
$ws.Range("B3").Characters(32, 5).Font.Bold = $true

# --- Notes for rows 3-5 ---
$ws.Range("D3").Value = "Cheaper than Mouser order, `$18 + `$7 shipping"
$ws.Range("D5").Value = "Saved `$15 on shipping with a first order coupon"
$ws.Range("D4").Value = "Went with FedEx International for about `$10 less shipping"

# --- Row 9: Total ---
$ws.Range("B9").Value = "Total:"
$ws.Range("B9").HorizontalAlignment = -4152
$ws.Range("C9").Formula = "=SUM(C2:C7)"

# --- Row 6: Amazon / SSD1306 128x64 OLED Display x 6 (bold " x 6") ---
$ws.Range("A6").Value = "Amazon"
$ws.Range("B6").Value = "SSD1306 128x64 OLED Display x 6"
$ws.Range("B6").Characters(28, 4).Font.Bold = $true
$ws.Range("D6").Value = "Bought 4 blue PCBs and 2 black PCBs, both white pixels"

# --- Row 7: Rotary encoder knobs (not ordered yet, no vendor/cost) ---
$ws.Range("B7").Value = "Rotary encoder knobs"

# --- Costs ---
$ws.Range("C4").Value = 34.91
$ws.Range("C5").Value = 59.85
$ws.Range("C6").Value = 33.32

# --- Column widths ---
$ws.Columns("D").ColumnWidth = 65.05338541666667

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- Selection ---
[void]$ws.Range("C10").Select()

Write-Host "done"
